# Generate Report for Handback
# Update the "Correspond Handback DateTime" values recorded for the
# zh-cn and de-de handback report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-02-22 09:24:19"
$wsZhCn.Range("G3").Value = "2016-02-22 09:25:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-02-22 09:24:31"
$wsDeDe.Range("G3").Value = "2016-02-22 09:25:25"
